# Apply the daily cryptos.xlsx price/volume refresh (GitHub Actions data pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '29.283.01'
$ws.Range('E2').Value = '  +0.22%  '

# Row 3
$ws.Range('D3').Value = '1.870.77'
$ws.Range('E3').Value = '  +0.20%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  +0.08%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7075'
$ws.Range('E5').Value = '  -0.57%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '241.58'

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.000'
$ws.Range('E7').Value = '  +0.05%  '

# Row 8
$ws.Range('B8').Value = 'Dogecoin'
$ws.Range('C8').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07761'
$ws.Range('E8').Value = '  +1.14%  '

# Row 9
$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3090'
$ws.Range('E9').Value = '  -0.86%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '25.02'
$ws.Range('E10').Value = '  +1.04%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08393'
$ws.Range('E11').Value = '  +0.32%  '

# Row 12
$ws.Range('D12').Value = '1.869.48'
$ws.Range('E12').Value = '  +0.30%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.237'
$ws.Range('E13').Value = '  +0.15%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.7111'
$ws.Range('E14').Value = '  -0.15%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '90.98'
$ws.Range('E15').Value = '  -0.45%  '

# Row 16
$ws.Range('D16').Value = '29.298.79'
$ws.Range('E16').Value = '  +0.25%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.058'
$ws.Range('E17').Value = '  +1.93%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008171'
$ws.Range('E18').Value = '  +4.56%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '239.75'
$ws.Range('E19').Value = '  -1.63%  '

# Row 20
$ws.Range('E20').Value = '  +0.65%  '

# Row 21
$ws.Range('D21').Value = '2.118.59'
$ws.Range('E21').Value = '  +0.20%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('E22').Value = '  +0.16%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.745'
$ws.Range('E23').Value = '  -1.60%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.000'
$ws.Range('E24').Value = '  +0.09%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1583'
$ws.Range('E25').Value = '  -0.76%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '163.35'
$ws.Range('E26').Value = '  -0.01%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.014'
$ws.Range('E27').Value = '  +0.71%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.46'
$ws.Range('E28').Value = '  -0.19%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.506'
$ws.Range('E29').Value = '  +0.41%  '

# Row 30
$ws.Range('E30').Value = '  -0.11%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.291'
$ws.Range('E31').Value = '  -2.42%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.299'
$ws.Range('E32').Value = '  +0.93%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05314'
$ws.Range('E33').Value = '  +2.99%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.934'
$ws.Range('E34').Value = '  +1.01%  '

# Row 35
$ws.Range('E35').Value = '  +0.69%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7432'
$ws.Range('E36').Value = '  -6.82%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.698'
$ws.Range('E37').Value = '  +0.52%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01869'
$ws.Range('E38').Value = '  +0.87%  '

# Row 39
$ws.Range('D39').Value = '1.232.46'
$ws.Range('E39').Value = '  +5.75%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.728'
$ws.Range('E40').Value = '  +0.63%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.564'
$ws.Range('E41').Value = '  +4.28%  '

# Row 42
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '109.85'
$ws.Range('E42').Value = '  +6.59%  '

# Row 43
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8860'
$ws.Range('E43').Value = '  -1.27%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '72.39'
$ws.Range('E44').Value = '  -1.22%  '

# Row 45
$ws.Range('E45').Value = '  +0.11%  '

# Row 46
$ws.Range('D46').Value = '2.016.05'
$ws.Range('E46').Value = '  +0.21%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5193'
$ws.Range('E47').Value = '  +0.17%  '

# Row 48
$ws.Range('E48').Value = '  +0.65%  '

# Row 49
$ws.Range('E49').Value = '  +1.73%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.390'

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4304'
$ws.Range('E51').Value = '  +0.18%  '
